# Add new column 'Servised by' to Card2 by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# New header cell O1, styled like the other header cells (bold/centered/bordered)
$ws.Range("O1").Value = "Servised by"
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in N2:N13 with "nan" (previously blank) and add blank O2:O13 cells
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 14).Value = "nan"   # column N
    # Force the new O cell into existence (as an empty cell) without
    # introducing a new style, mirroring the other blank cells on the sheet.
    $ws.Cells.Item($row, 15).Borders.LineStyle = -4142   # xlLineStyleNone
}
